$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) for each crypto row.
# Leading apostrophe forces text interpretation so values like '1.00'/'39.70' keep
# trailing zeros and multi-dot strings like "51.734.87" are not parsed as numbers.
$ws.Range("D2").Value = "'51.734.87"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = "'2.768.45"
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = "'357.47"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = "'108.67"
$ws.Range("E6").Value = '  -4.86%  '
$ws.Range("D7").Value = "'0.556"
$ws.Range("E7").Value = '  +1.70%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = "'0.587"
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").Value = "'39.70"
$ws.Range("E10").Value = '  -5.02%  '
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = "'19.36"
$ws.Range("E13").Value = '  -2.98%  '
$ws.Range("D14").Value = "'7.59"
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = "'3.213.96"
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").Value = "'2.775.35"
$ws.Range("E16").Value = '  -3.46%  '
$ws.Range("D17").Value = "'0.909"
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = "'51.662.07"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = "'7.37"
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").Value = "'3.09"
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").Value = "'12.97"
$ws.Range("E21").Value = '  -5.82%  '
$ws.Range("D22").Value = "'0.0₃0973"
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = "'272.76"
$ws.Range("E23").Value = '  +1.47%  '
$ws.Range("D24").Value = "'69.28"
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("D25").Value = "'2.74"
$ws.Range("E25").Value = '  -2.43%  '
$ws.Range("D26").Value = "'26.32"
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").Value = "'0.0467"
$ws.Range("E31").Value = '  +6.05%  '
$ws.Range("D32").Value = "'51.12"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").Value = "'33.90"
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("D34").Value = "'5.70"
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("D35").Value = "'5.33"
$ws.Range("E35").Value = '  +7.88%  '
$ws.Range("D36").Value = "'0.0833"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("E39").Value = '  -5.71%  '
$ws.Range("D40").Value = "'17.93"
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D41").Value = "'0.114"
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("D42").Value = "'125.05"
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("D43").Value = "'2.50"
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("E44").Value = '  -1.91%  '
$ws.Range("D45").Value = "'21.83"
$ws.Range("E45").Value = '  -6.67%  '
$ws.Range("D46").Value = "'2.052.08"
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").Value = "'3.20"
$ws.Range("E48").Value = '  -4.79%  '
$ws.Range("D49").Value = "'5.66"
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").Value = "'0.922"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("D51").Value = "'8.94"
$ws.Range("E51").Value = '  -0.37%  '